# Add an "Address" column (new column F) between "After training"/"Name"
# columns and the existing "District" column, which shifts right from F to G.
# The address text is derived from the second line of the "Names" cell in
# column B, which holds "<Name>\n<Address>, <District>." - everything
# before the final comma-separated District segment becomes the Address
# (joined back together without the comma separators), and the trailing
# period is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; the old F ("District") becomes G.
$ws.Columns("F:F").Insert()

# Header for the newly inserted column.
$ws.Range("F2").Value = "Address"

# Sheet has data rows 3 through 55 (row 1 = title, row 2 = headers).
$firstRow = 3
$lastRow = 55

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $nameCell = $ws.Cells.Item($r, 2)   # column B
    $raw = $nameCell.Value()
    if ($raw -eq $null) { continue }

    $lines = $raw -split "`n"
    if ($lines.Count -lt 2) { continue }

    $addressLine = $lines[1].TrimEnd(".").Trim()
    $segments = $addressLine -split ","

    if ($segments.Count -lt 1) { continue }

    # Every comma-separated segment except the last (the District) makes
    # up the Address; they are concatenated back together with no
    # separator, matching the source data.
    $addressSegments = $segments[0..($segments.Count - 2)]

    $address = ""
    foreach ($seg in $addressSegments) {
        $address = $address + $seg.Trim()
    }

    $ws.Cells.Item($r, 6).Value = $address   # column F
}
